# Commit: Fixed update to excel issue
# - Rename "Requested quantity" headers to per-sheet metric names
# - Add new "PO Forecast" sheet with ds/PO_Forecast/yhat_lower/yhat_upper forecast data

$wb = $excel.ActiveWorkbook

# --- 1) Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after "Monthly Trend" ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3) Match formatting used on the other sheets: bold/bordered header row (A1:D1)
#        and the date number format on column A (rows 2-47) ---
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A47").PasteSpecial(-4122)

$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# --- 4) Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 5) Forecast data rows 2-47 (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$wsForecast.Cells.Item(2, 1).Value = 45165.99999999999
$wsForecast.Cells.Item(2, 2).Value = 93
$wsForecast.Cells.Item(2, 3).Value = -161.0832541714354
$wsForecast.Cells.Item(2, 4).Value = 352.2386454247885
$wsForecast.Cells.Item(3, 1).Value = 45186.99999999999
$wsForecast.Cells.Item(3, 2).Value = 102
$wsForecast.Cells.Item(3, 3).Value = -155.4476065709814
$wsForecast.Cells.Item(3, 4).Value = 359.8139361912179
$wsForecast.Cells.Item(4, 1).Value = 45207.99999999999
$wsForecast.Cells.Item(4, 2).Value = 111
$wsForecast.Cells.Item(4, 3).Value = -148.125759078578
$wsForecast.Cells.Item(4, 4).Value = 381.3984711386434
$wsForecast.Cells.Item(5, 1).Value = 45214.99999999999
$wsForecast.Cells.Item(5, 2).Value = 115
$wsForecast.Cells.Item(5, 3).Value = -158.0411168715631
$wsForecast.Cells.Item(5, 4).Value = 368.9129765977451
$wsForecast.Cells.Item(6, 1).Value = 45221.99999999999
$wsForecast.Cells.Item(6, 2).Value = 118
$wsForecast.Cells.Item(6, 3).Value = -167.338456843692
$wsForecast.Cells.Item(6, 4).Value = 369.9472648611818
$wsForecast.Cells.Item(7, 1).Value = 45235.99999999999
$wsForecast.Cells.Item(7, 2).Value = 124
$wsForecast.Cells.Item(7, 3).Value = -121.5467610084868
$wsForecast.Cells.Item(7, 4).Value = 377.6330051077099
$wsForecast.Cells.Item(8, 1).Value = 45242.99999999999
$wsForecast.Cells.Item(8, 2).Value = 127
$wsForecast.Cells.Item(8, 3).Value = -134.9415875754592
$wsForecast.Cells.Item(8, 4).Value = 392.6157215238615
$wsForecast.Cells.Item(9, 1).Value = 45249.99999999999
$wsForecast.Cells.Item(9, 2).Value = 130
$wsForecast.Cells.Item(9, 3).Value = -129.7656068545577
$wsForecast.Cells.Item(9, 4).Value = 409.1532189692283
$wsForecast.Cells.Item(10, 1).Value = 45256.99999999999
$wsForecast.Cells.Item(10, 2).Value = 133
$wsForecast.Cells.Item(10, 3).Value = -119.1444838774471
$wsForecast.Cells.Item(10, 4).Value = 406.5758677543001
$wsForecast.Cells.Item(11, 1).Value = 45263.99999999999
$wsForecast.Cells.Item(11, 2).Value = 136
$wsForecast.Cells.Item(11, 3).Value = -133.8693017092854
$wsForecast.Cells.Item(11, 4).Value = 395.7346476611111
$wsForecast.Cells.Item(12, 1).Value = 45270.99999999999
$wsForecast.Cells.Item(12, 2).Value = 139
$wsForecast.Cells.Item(12, 3).Value = -131.4449947882444
$wsForecast.Cells.Item(12, 4).Value = 398.9002052044934
$wsForecast.Cells.Item(13, 1).Value = 45277.99999999999
$wsForecast.Cells.Item(13, 2).Value = 142
$wsForecast.Cells.Item(13, 3).Value = -136.8967522049984
$wsForecast.Cells.Item(13, 4).Value = 411.6954364933451
$wsForecast.Cells.Item(14, 1).Value = 45298.99999999999
$wsForecast.Cells.Item(14, 2).Value = 151
$wsForecast.Cells.Item(14, 3).Value = -110.2373018402158
$wsForecast.Cells.Item(14, 4).Value = 400.7039802698002
$wsForecast.Cells.Item(15, 1).Value = 45305.99999999999
$wsForecast.Cells.Item(15, 2).Value = 154
$wsForecast.Cells.Item(15, 3).Value = -109.7183059018758
$wsForecast.Cells.Item(15, 4).Value = 412.4855947672582
$wsForecast.Cells.Item(16, 1).Value = 45312.99999999999
$wsForecast.Cells.Item(16, 2).Value = 158
$wsForecast.Cells.Item(16, 3).Value = -122.8069289675928
$wsForecast.Cells.Item(16, 4).Value = 412.4276235866313
$wsForecast.Cells.Item(17, 1).Value = 45326.99999999999
$wsForecast.Cells.Item(17, 2).Value = 164
$wsForecast.Cells.Item(17, 3).Value = -103.635139640505
$wsForecast.Cells.Item(17, 4).Value = 457.4223078809387
$wsForecast.Cells.Item(18, 1).Value = 45333.99999999999
$wsForecast.Cells.Item(18, 2).Value = 167
$wsForecast.Cells.Item(18, 3).Value = -101.7662388071602
$wsForecast.Cells.Item(18, 4).Value = 434.7002322222405
$wsForecast.Cells.Item(19, 1).Value = 45347.99999999999
$wsForecast.Cells.Item(19, 2).Value = 173
$wsForecast.Cells.Item(19, 3).Value = -87.74576863858704
$wsForecast.Cells.Item(19, 4).Value = 433.321147785381
$wsForecast.Cells.Item(20, 1).Value = 45361.99999999999
$wsForecast.Cells.Item(20, 2).Value = 179
$wsForecast.Cells.Item(20, 3).Value = -105.7268301545154
$wsForecast.Cells.Item(20, 4).Value = 457.5075624442862
$wsForecast.Cells.Item(21, 1).Value = 45368.99999999999
$wsForecast.Cells.Item(21, 2).Value = 182
$wsForecast.Cells.Item(21, 3).Value = -81.48240147034102
$wsForecast.Cells.Item(21, 4).Value = 471.0061466701038
$wsForecast.Cells.Item(22, 1).Value = 45375.99999999999
$wsForecast.Cells.Item(22, 2).Value = 185
$wsForecast.Cells.Item(22, 3).Value = -81.58119115554213
$wsForecast.Cells.Item(22, 4).Value = 464.5118902128983
$wsForecast.Cells.Item(23, 1).Value = 45382.99999999999
$wsForecast.Cells.Item(23, 2).Value = 188
$wsForecast.Cells.Item(23, 3).Value = -71.81217512042819
$wsForecast.Cells.Item(23, 4).Value = 480.5196784444874
$wsForecast.Cells.Item(24, 1).Value = 45403.99999999999
$wsForecast.Cells.Item(24, 2).Value = 197
$wsForecast.Cells.Item(24, 3).Value = -86.15531796950233
$wsForecast.Cells.Item(24, 4).Value = 474.4338570465708
$wsForecast.Cells.Item(25, 1).Value = 45417.99999999999
$wsForecast.Cells.Item(25, 2).Value = 204
$wsForecast.Cells.Item(25, 3).Value = -65.23457032270302
$wsForecast.Cells.Item(25, 4).Value = 465.7464331009515
$wsForecast.Cells.Item(26, 1).Value = 45424.99999999999
$wsForecast.Cells.Item(26, 2).Value = 207
$wsForecast.Cells.Item(26, 3).Value = -45.79987593244207
$wsForecast.Cells.Item(26, 4).Value = 474.4269004232962
$wsForecast.Cells.Item(27, 1).Value = 45431.99999999999
$wsForecast.Cells.Item(27, 2).Value = 210
$wsForecast.Cells.Item(27, 3).Value = -47.67359284613406
$wsForecast.Cells.Item(27, 4).Value = 477.4711941652171
$wsForecast.Cells.Item(28, 1).Value = 45445.99999999999
$wsForecast.Cells.Item(28, 2).Value = 216
$wsForecast.Cells.Item(28, 3).Value = -62.94359808334428
$wsForecast.Cells.Item(28, 4).Value = 475.1373695155643
$wsForecast.Cells.Item(29, 1).Value = 45459.99999999999
$wsForecast.Cells.Item(29, 2).Value = 222
$wsForecast.Cells.Item(29, 3).Value = -40.97806386144683
$wsForecast.Cells.Item(29, 4).Value = 512.3021327346605
$wsForecast.Cells.Item(30, 1).Value = 45466.99999999999
$wsForecast.Cells.Item(30, 2).Value = 225
$wsForecast.Cells.Item(30, 3).Value = -50.31006005744987
$wsForecast.Cells.Item(30, 4).Value = 478.882481181762
$wsForecast.Cells.Item(31, 1).Value = 45473.99999999999
$wsForecast.Cells.Item(31, 2).Value = 228
$wsForecast.Cells.Item(31, 3).Value = -41.69581632371026
$wsForecast.Cells.Item(31, 4).Value = 472.9497759212175
$wsForecast.Cells.Item(32, 1).Value = 45515.99999999999
$wsForecast.Cells.Item(32, 2).Value = 247
$wsForecast.Cells.Item(32, 3).Value = -25.81042266611511
$wsForecast.Cells.Item(32, 4).Value = 518.6382687883157
$wsForecast.Cells.Item(33, 1).Value = 45543.99999999999
$wsForecast.Cells.Item(33, 2).Value = 259
$wsForecast.Cells.Item(33, 3).Value = -8.810552700951169
$wsForecast.Cells.Item(33, 4).Value = 534.4969419346054
$wsForecast.Cells.Item(34, 1).Value = 45550.99999999999
$wsForecast.Cells.Item(34, 2).Value = 262
$wsForecast.Cells.Item(34, 3).Value = -19.98170383241752
$wsForecast.Cells.Item(34, 4).Value = 524.2766183278387
$wsForecast.Cells.Item(35, 1).Value = 45557.99999999999
$wsForecast.Cells.Item(35, 2).Value = 265
$wsForecast.Cells.Item(35, 3).Value = -13.47892453880151
$wsForecast.Cells.Item(35, 4).Value = 532.2750426413132
$wsForecast.Cells.Item(36, 1).Value = 45578.99999999999
$wsForecast.Cells.Item(36, 2).Value = 274
$wsForecast.Cells.Item(36, 3).Value = 6.3821048472648
$wsForecast.Cells.Item(36, 4).Value = 539.1392102984345
$wsForecast.Cells.Item(37, 1).Value = 45585.99999999999
$wsForecast.Cells.Item(37, 2).Value = 277
$wsForecast.Cells.Item(37, 3).Value = 14.876259362627
$wsForecast.Cells.Item(37, 4).Value = 558.8467301856571
$wsForecast.Cells.Item(38, 1).Value = 45599.99999999999
$wsForecast.Cells.Item(38, 2).Value = 283
$wsForecast.Cells.Item(38, 3).Value = 12.45113846400483
$wsForecast.Cells.Item(38, 4).Value = 532.7081232111134
$wsForecast.Cells.Item(39, 1).Value = 45627.99999999999
$wsForecast.Cells.Item(39, 2).Value = 296
$wsForecast.Cells.Item(39, 3).Value = 30.00785645409023
$wsForecast.Cells.Item(39, 4).Value = 570.6052388663635
$wsForecast.Cells.Item(40, 1).Value = 45634.99999999999
$wsForecast.Cells.Item(40, 2).Value = 299
$wsForecast.Cells.Item(40, 3).Value = 37.73722034356455
$wsForecast.Cells.Item(40, 4).Value = 554.2219731692589
$wsForecast.Cells.Item(41, 1).Value = 45641.99999999999
$wsForecast.Cells.Item(41, 2).Value = 302
$wsForecast.Cells.Item(41, 3).Value = 35.10693610603481
$wsForecast.Cells.Item(41, 4).Value = 601.905771041164
$wsForecast.Cells.Item(42, 1).Value = 45648.99999999999
$wsForecast.Cells.Item(42, 2).Value = 305
$wsForecast.Cells.Item(42, 3).Value = 38.85982280746453
$wsForecast.Cells.Item(42, 4).Value = 555.1329506036952
$wsForecast.Cells.Item(43, 1).Value = 45655.99999999999
$wsForecast.Cells.Item(43, 2).Value = 308
$wsForecast.Cells.Item(43, 3).Value = 52.92097838976629
$wsForecast.Cells.Item(43, 4).Value = 560.8324218768325
$wsForecast.Cells.Item(44, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(44, 2).Value = 311
$wsForecast.Cells.Item(44, 3).Value = 26.99789894843342
$wsForecast.Cells.Item(44, 4).Value = 577.5506898189988
$wsForecast.Cells.Item(45, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(45, 2).Value = 314
$wsForecast.Cells.Item(45, 3).Value = 25.74751556694002
$wsForecast.Cells.Item(45, 4).Value = 578.5414950090673
$wsForecast.Cells.Item(46, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(46, 2).Value = 317
$wsForecast.Cells.Item(46, 3).Value = 25.09610181546724
$wsForecast.Cells.Item(46, 4).Value = 585.3455508873226
$wsForecast.Cells.Item(47, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(47, 2).Value = 320
$wsForecast.Cells.Item(47, 3).Value = 47.56727876576107
$wsForecast.Cells.Item(47, 4).Value = 567.4775972013646

Write-Output "Done: renamed headers + added PO Forecast sheet ($($wsForecast.UsedRange.Rows.Count) rows)"
